# Generate Report for Handback
# Update the timestamp cells recorded in the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Sheets.Item("Overview")
$wsZhCn     = $wb.Sheets.Item("zh-cn")
$wsDeDe     = $wb.Sheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file
$wsOverview.Range("G2").Value = "2016-08-27 11:04:57"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for first row
$wsZhCn.Range("H2").Value = "2016-08-27 11:04:53"
$wsZhCn.Range("K2").Value = "2016-08-27 11:05:16"

# de-de sheet: Correspond Handoff Datetime (same value as Overview) / Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-08-27 11:04:57"
$wsDeDe.Range("K2").Value = "2016-08-27 11:05:23"
